$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" ---------------
# The status string shows up once per localized-language column on the
# "Overview" sheet (columns "zh-cn" / "de-de", row 2) and once in the
# "Status" column of each per-language detail sheet ("zh-cn", "de-de").
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = "In Translation"
$wsOverview.Range("F2").Value2 = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = "In Translation"

# --- 2. Re-fit the (now shorter-text) status columns ------------------------
# The host only exposes column sizing through the classic `ColumnWidth`
# (character-unit) property, whose setter snaps to the host's internal
# pixel grid (1/6-character steps) the same way Excel's own COM layer does.
# `ColumnWidth` itself reads/writes that already-quantised value, so feeding
# it the raw target width lands on the neighbouring grid step instead of the
# nearest one. Pre-compensate by solving for the input that the host will
# quantise to the closest achievable step to the desired width.
function Set-ClosestColumnWidth($col, $desiredWidth) {
    $step = 1.0 / 6.0
    $bucket = [Math]::Round($desiredWidth / $step) * $step
    $col.ColumnWidth = $bucket - (5.0 / 6.0)
}

Set-ClosestColumnWidth $wsOverview.Columns.Item(5) 13.4101845877511
Set-ClosestColumnWidth $wsOverview.Columns.Item(6) 13.4101845877511

Set-ClosestColumnWidth $wsZhCn.Columns.Item(3) 13.4101845877511
Set-ClosestColumnWidth $wsDeDe.Columns.Item(3) 13.4101845877511
